# Add 5 new Solar Power Station producer rows to the "Electric Network" sheet.
# Commit message: "Added additional producers (+5 solar power stations)"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Electric Network
$ws2 = $wb.Worksheets.Item(2)   # Gas Network

# ---------------------------------------------------------------------------
# 1. Build the 5 new producer rows on "Electric Network" (rows 7-11), which
#    previously only held the stray "Transmitter"/"Consumer" labels in B7/B8.
# ---------------------------------------------------------------------------

$newRows = @(
    @{ Row=7;  Name="Liberty 1";            Power=259; Loc="29.9589, -94.8505";  East=321434.38299999997; North=3315671.4759999998; Year=2024 },
    @{ Row=8;  Name="Trinity river solar";  Power=409; Loc="30.0103, -94.4863 "; East=356656.51699999999; North=3320856.71;         Year=2025 },
    @{ Row=9;  Name="Myrtle solar";         Power=538; Loc="29.2285, -95.4274";  East=264056.06400000001; North=3235743.32;         Year=2023 },
    @{ Row=10; Name="Red Bluff Road Solar"; Power=360; Loc="29.6150, -95.0702";  East=299545.32500000001; North=3277915.5559999999; Year=2100 },
    @{ Row=11; Name="Brazoria West";        Power=540; Loc="29.1910, -95.6630";  East=241053.016;         North=3232083.39;         Year=2022 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy the formatting of the shaded producer row 4 (B4:I4) onto the new
    # row so the fills / alignments / font used for "Producer" rows match.
    $ws1.Range("B4:I4").Copy()
    $destRow = $ws1.Range($ws1.Cells.Item($row, 2), $ws1.Cells.Item($row, 9))
    $destRow.PasteSpecial(-4122)   # xlPasteFormats

    # "Asset type" (column D) is never shaded on any existing producer row,
    # so strip its pasted formatting back to the default style.
    $ws1.Cells.Item($row, 4).ClearFormats()

    $ws1.Cells.Item($row, 2).Value = "Producer "
    $ws1.Cells.Item($row, 3).Value = $r.Name
    $ws1.Cells.Item($row, 4).Value = "Solar Power Station"
    $ws1.Cells.Item($row, 5).Value = $r.Power
    $ws1.Cells.Item($row, 6).Value = $r.Loc
    $ws1.Cells.Item($row, 7).Value = $r.East
    $ws1.Cells.Item($row, 8).Value = $r.North

    # Easting/Northing (G/H) on the new rows use a plain "General" number
    # format instead of the shaded style copied above.
    $destGH = $ws1.Range($ws1.Cells.Item($row, 7), $ws1.Cells.Item($row, 8))
    $destGH.NumberFormat = "General"

    $ws1.Cells.Item($row, 9).Value = $r.Year
}

# ---------------------------------------------------------------------------
# 2. Update sheet selections / active tab to match the saved view state.
#    "Electric Network" becomes the selected tab, "Gas Network" selection
#    moves to E15.
# ---------------------------------------------------------------------------

$ws2.Activate()
$ws2.Range("E15").Select()

$ws1.Activate()
$ws1.Range("E13").Select()
